# Task #1094: Multiple payees
# Update the Payee column (column B) on the "Definition" sheet so a handful
# of rows reference multiple payees instead of a single one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definition")

# Housing:Home Goods -> payee now lists several stores
$ws.Range("B14").Value = "Bed Bath & Beyond,Target,Container Store"

# Income:Bonus -> payee now "Megacorp, Inc."
$ws.Range("B45").Value = "Megacorp, Inc."

# Food:Away:Dinner -> payee now lists several restaurants
$ws.Range("B17").Value = "Applebees,Olive Garden,Spaghetti Factory"

# Food:Away:Coffee -> payee now lists several coffee shops
$ws.Range("B16").Value = "Starbucks,Uptown Espresso,Tim Horton's"

# Income:Salary -> payee now "Megacorp Inc."
$ws.Range("B33").Value = "Megacorp Inc."

# Move the visible selection/active cell to B34, matching the saved view.
$ws.Range("B34").Select() | Out-Null
